$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 360 so the existing Jengibre price history
# (currently rows 360-387) shifts down to rows 361-388, and fill the new
# row 360 with the latest weekly price entry.
$ws.Rows.Item(360).Insert()

$ws.Range("A360").Value = 10
$ws.Range("B360").Value = "Vega Modelo de Temuco"
$ws.Range("C360").Value = "La Araucanía"
$ws.Range("D360").Value = 45265
$ws.Range("D360").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E360").Value = 9
$ws.Range("F360").Value = 100114007
$ws.Range("G360").Value = "Jengibre"
$ws.Range("H360").Value = "Sin especificar"
$ws.Range("I360").Value = "Primera"
$ws.Range("J360").Value = 100
$ws.Range("K360").Value = 26000
$ws.Range("L360").Value = 26000
$ws.Range("M360").Value = 26000
$ws.Range("N360").Value = "$/caja 13 kilos"
$ws.Range("O360").Value = "Perú"
$ws.Range("P360").Value = 2000
$ws.Range("Q360").Value = 13
$ws.Range("R360").Value = "Hortaliza"
